$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cars slow down sooner / accelerate within safe limits: shrink the
# "Safe Accel. Distance" curve in row 9 (AF:AJ get new, lower values) and
# the vehicle now reaches top speed by AJ9, so AK9:AT9 no longer apply.
$ws.Range("AF9").Value = 286
$ws.Range("AG9").Value = 294
$ws.Range("AH9").Value = 300
$ws.Range("AI9").Value = 304
$ws.Range("AJ9").Value = 306
$ws.Range("AK9:AT9").ClearContents()

# Recalculate so the row 11 "host safe dist." formulas (SUM(x9-x18)) pick
# up the new row 9 values / now-blank AK9:AT9 cells.
$excel.Calculate()

# The tailing-distance column AK:AT now displays negative 3-digit values
# (e.g. -327), one character wider than before, so Excel's best-fit
# widens those columns versus AD:AJ.
$ws.Columns("AK:AT").ColumnWidth = 3.5

# Leave the selection where the author ended up reviewing the new numbers.
$ws.Range("AU11").Select()
